# Fix FEKT registration link.
#
# The "FEKT" bullet used to read:   FEKT - goo.gl/cBXSLd   (hyphen, link split
# across two runs "goo.gl/" + "cBXSLd", hyperlink flagged with the stale
# ppaction://hlinkfile action).
#
# It should read:                   FEKT – http://goo.gl/cBXSLd
# with a single run holding the (now schemed) URL and a plain hyperlink.

$p = $ppt.ActivePresentation

$needle        = "FEKT - goo.gl/cBXSLd"
$prefixOld     = "FEKT - "
$prefixNew     = "FEKT " + [char]0x2013 + " "   # en dash
$linkTextOld   = "goo.gl/cBXSLd"
$linkTextNew   = "http://goo.gl/cBXSLd"

# Locate the shape holding the text, rather than hard-coding a slide index.
$targetSlide = $null
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $shp = $sl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t.IndexOf($needle) -ge 0) {
                    $targetSlide = $sl
                    $targetShape = $shp
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text
    $idx0 = $full.IndexOf($needle)
    $startPos = $idx0 + 1   # TextRange.Characters is 1-based

    # 1) "FEKT - " -> "FEKT \u2013 " (hyphen becomes en dash)
    $prefixRange = $tr.Characters($startPos, $prefixOld.Length)
    $prefixRange.Text = $prefixNew

    # 2) Rewrite the link text to add the http:// scheme. Setting .Text across
    #    the old two-run span ("goo.gl/" + "cBXSLd") also merges them into one
    #    run, matching the target markup.
    $linkStartPos = $startPos + $prefixOld.Length
    $linkRange = $tr.Characters($linkStartPos, $linkTextOld.Length)
    $linkRange.Text = $linkTextNew

    # 3) Re-apply the (unchanged) hyperlink target so the host regenerates a
    #    clean <a:hlinkClick r:id="..."/> without the stale
    #    action="ppaction://hlinkfile" marker.
    $linkRange2 = $tr.Characters($linkStartPos, $linkTextNew.Length)
    $hyperlink = $linkRange2.ActionSettings(1).Hyperlink
    $hyperlink.Address = $linkTextOld
}
